$d = $word.ActiveDocument

# 1) "为消费者提供退换货后的包运费服务，并收取商家服务费。"
#    -> "为消费者提供退换货后的包运费投保、理赔服务，并收取商家相应的服务费。"
$d.Content.Find.Execute(
    "为消费者提供退换货后的包运费服务，并收取商家服务费。", $true, $false, $false, $false, $false,
    $true, 1, $false, "为消费者提供退换货后的包运费投保、理赔服务，并收取商家相应的服务费。", 2) | Out-Null

# 2) " 通过" -> " 通过离线"  (the run right before the "Hive" split text)
$d.Content.Find.Execute(
    " 通过", $true, $false, $false, $false, $false,
    $true, 1, $false, " 通过离线", 2) | Out-Null

# 3) "ve离线指标加工成的风控策略与风险定价，来控制商家的金融风险与每日服务费定价"
#    -> "ve数据生成的指标加工成的风控策略与风险定价，来控制商家的金融风险与每日服务费定价"
$d.Content.Find.Execute(
    "ve离线指标加工成的风控策略与风险定价，来控制商家的金融风险与每日服务费定价", $true, $false, $false, $false, $false,
    $true, 1, $false, "ve数据生成的指标加工成的风控策略与风险定价，来控制商家的金融风险与每日服务费定价", 2) | Out-Null

# 4) "在模型与数据抽象的基础上，通过权益配置，来快速支持开店礼包、升级版等业务需求"
#    -> "在模型与数据抽象的基础上，通过权益配置，来快速支持开店礼包、升级版服务等业务需求"
$d.Content.Find.Execute(
    "在模型与数据抽象的基础上，通过权益配置，来快速支持开店礼包、升级版等业务需求", $true, $false, $false, $false, $false,
    $true, 1, $false, "在模型与数据抽象的基础上，通过权益配置，来快速支持开店礼包、升级版服务等业务需求", 2) | Out-Null

# 5) "  1. 有赞寄件是在交易支付后物流下单、结算的核心，通过物流商运营、运费定价等构造整个物流计费结算模型"
#    -> "  1. 有赞寄件是在交易支付后物流下单、结算的核心，其通过物流商运营、运费定价等构造整个物流计费结算模型"
$d.Content.Find.Execute(
    "  1. 有赞寄件是在交易支付后物流下单、结算的核心，通过物流商运营、运费定价等构造整个物流计费结算模型", $true, $false, $false, $false, $false,
    $true, 1, $false, "  1. 有赞寄件是在交易支付后物流下单、结算的核心，其通过物流商运营、运费定价等构造整个物流计费结算模型", 2) | Out-Null

# 6) "2. 上门取件是作为交易逆向售后的核心服务，通过状态机+？保证上门取件单、三方物流单、交易单之间状态一致"
#    -> "2. 上门取件是作为交易逆向售后的核心服务，并通过状态机+？保证上门取件单、三方物流单、交易单之间状态一致"
$d.Content.Find.Execute(
    "2. 上门取件是作为交易逆向售后的核心服务，通过状态机+？保证上门取件单、三方物流单、交易单之间状态一致", $true, $false, $false, $false, $false,
    $true, 1, $false, "2. 上门取件是作为交易逆向售后的核心服务，并通过状态机+？保证上门取件单、三方物流单、交易单之间状态一致", 2) | Out-Null

# 7) Move the "_GoBack" last-edit-position bookmark to sit right after the
#    "退货包运费" heading run (matching where Word last left the cursor after
#    this editing session).
$headingRng = $d.Content
$headingRng.Find.Execute("退货包运费") | Out-Null
$headingRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $headingRng) | Out-Null
